$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" -------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("D15").Value = 1267.19
$ws1.Range("D19").Value = "1 de 17"

# --- Sheet "VENTA MENSUAL" -----------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F15").Value = 1267.19
$ws2.Range("F19").Value = 1362.45
# column F widened slightly (12 -> 13 raw OOXML units) to fit the new value
$ws2.Columns("F").ColumnWidth = 12.166666666666666

# --- Sheet "CUMPLIMIENTO MENSUAL" ----------------------------------------
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D3").Value = 1267.19
$ws3.Range("E3").Value = 7401.719999999999
$ws3.Range("F3").Value = 0.1461763935719716
$ws3.Range("D19").Value = 1362.45
$ws3.Range("E19").Value = 45856.85386304604
$ws3.Range("F19").Value = 0.02885366552526111
# columns D and E widened slightly (11 -> 13, 22 -> 23 raw OOXML units)
$ws3.Columns("D").ColumnWidth = 12.166666666666666
$ws3.Columns("E").ColumnWidth = 22.166666666666668
